$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new training-session column AD for 2025-08-19 (serial 45888) ---

# Header cell AD1: copy date-cell format from J1 (same style as other "early" date headers, s=3)
$ws.Range("J1").Copy() | Out-Null
$ws.Range("AD1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("AD1").Value2 = 45888

# Data cells AD2:AD27: copy format from AC2:AC27 (s=4)
$ws.Range("AC2:AC27").Copy() | Out-Null
$ws.Range("AD2:AD27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$values = @("P","B","P","P","A","P","P","P","P","P","P","A","P","P","P","B","P","RH","P","P","P","P","P","P","P","RH")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 30).Value2 = $values[$i]
}

# --- Remove the old summary row 28 (COUNTIF(AC2:AC27,"P") row no longer needed) ---
$ws.Rows.Item(28).Delete() | Out-Null

# Recalculate all formulas so cached <v> values reflect the new column
$excel.Calculate()

# Update the active selection to match the new cursor position
$ws.Range("AF25").Select() | Out-Null

$excel.CutCopyMode = 0
